$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-93 with new timestamp (2 days later) and new load values
$ws.Range("A2").Value = 45424
$ws.Range("B2").Value = 5000
$ws.Range("A3").Value = 45424.01041666666
$ws.Range("B3").Value = 4950
$ws.Range("A4").Value = 45424.02083333334
$ws.Range("B4").Value = 4900
$ws.Range("A5").Value = 45424.03125
$ws.Range("B5").Value = 4850
$ws.Range("A6").Value = 45424.04166666666
$ws.Range("B6").Value = 4800
$ws.Range("A7").Value = 45424.05208333334
$ws.Range("B7").Value = 4760
$ws.Range("A8").Value = 45424.0625
$ws.Range("B8").Value = 4730
$ws.Range("A9").Value = 45424.07291666666
$ws.Range("B9").Value = 4710
$ws.Range("A10").Value = 45424.08333333334
$ws.Range("B10").Value = 4700
$ws.Range("A11").Value = 45424.09375
$ws.Range("B11").Value = 4690
$ws.Range("A12").Value = 45424.10416666666
$ws.Range("B12").Value = 4680
$ws.Range("A13").Value = 45424.11458333334
$ws.Range("B13").Value = 4680
$ws.Range("A14").Value = 45424.125
$ws.Range("B14").Value = 4680
$ws.Range("A15").Value = 45424.13541666666
$ws.Range("B15").Value = 4680
$ws.Range("A16").Value = 45424.14583333334
$ws.Range("B16").Value = 4670
$ws.Range("A17").Value = 45424.15625
$ws.Range("B17").Value = 4670
$ws.Range("A18").Value = 45424.16666666666
$ws.Range("B18").Value = 4670
$ws.Range("A19").Value = 45424.17708333334
$ws.Range("B19").Value = 4670
$ws.Range("A20").Value = 45424.1875
$ws.Range("B20").Value = 4670
$ws.Range("A21").Value = 45424.19791666666
$ws.Range("B21").Value = 4670
$ws.Range("A22").Value = 45424.20833333334
$ws.Range("B22").Value = 4670
$ws.Range("A23").Value = 45424.21875
$ws.Range("B23").Value = 4680
$ws.Range("A24").Value = 45424.22916666666
$ws.Range("B24").Value = 4690
$ws.Range("A25").Value = 45424.23958333334
$ws.Range("B25").Value = 4700
$ws.Range("A26").Value = 45424.25
$ws.Range("B26").Value = 4720
$ws.Range("A27").Value = 45424.26041666666
$ws.Range("B27").Value = 4750
$ws.Range("A28").Value = 45424.27083333334
$ws.Range("B28").Value = 4770
$ws.Range("A29").Value = 45424.28125
$ws.Range("B29").Value = 4800
$ws.Range("A30").Value = 45424.29166666666
$ws.Range("B30").Value = 4820
$ws.Range("A31").Value = 45424.30208333334
$ws.Range("B31").Value = 4840
$ws.Range("A32").Value = 45424.3125
$ws.Range("B32").Value = 4850
$ws.Range("A33").Value = 45424.32291666666
$ws.Range("B33").Value = 4850
$ws.Range("A34").Value = 45424.33333333334
$ws.Range("B34").Value = 4830
$ws.Range("A35").Value = 45424.34375
$ws.Range("B35").Value = 4790
$ws.Range("A36").Value = 45424.35416666666
$ws.Range("B36").Value = 4750
$ws.Range("A37").Value = 45424.36458333334
$ws.Range("B37").Value = 4700
$ws.Range("A38").Value = 45424.375
$ws.Range("B38").Value = 4640
$ws.Range("A39").Value = 45424.38541666666
$ws.Range("B39").Value = 4590
$ws.Range("A40").Value = 45424.39583333334
$ws.Range("B40").Value = 4540
$ws.Range("A41").Value = 45424.40625
$ws.Range("B41").Value = 4490
$ws.Range("A42").Value = 45424.41666666666
$ws.Range("B42").Value = 4460
$ws.Range("A43").Value = 45424.42708333334
$ws.Range("B43").Value = 4440
$ws.Range("A44").Value = 45424.4375
$ws.Range("B44").Value = 4420
$ws.Range("A45").Value = 45424.44791666666
$ws.Range("B45").Value = 4410
$ws.Range("A46").Value = 45424.45833333334
$ws.Range("B46").Value = 4410
$ws.Range("A47").Value = 45424.46875
$ws.Range("B47").Value = 4400
$ws.Range("A48").Value = 45424.47916666666
$ws.Range("B48").Value = 4390
$ws.Range("A49").Value = 45424.48958333334
$ws.Range("B49").Value = 4380
$ws.Range("A50").Value = 45424.5
$ws.Range("B50").Value = 4360
$ws.Range("A51").Value = 45424.51041666666
$ws.Range("B51").Value = 4350
$ws.Range("A52").Value = 45424.52083333334
$ws.Range("B52").Value = 4340
$ws.Range("A53").Value = 45424.53125
$ws.Range("B53").Value = 4330
$ws.Range("A54").Value = 45424.54166666666
$ws.Range("B54").Value = 4330
$ws.Range("A55").Value = 45424.55208333334
$ws.Range("B55").Value = 4330
$ws.Range("A56").Value = 45424.5625
$ws.Range("B56").Value = 4340
$ws.Range("A57").Value = 45424.57291666666
$ws.Range("B57").Value = 4340
$ws.Range("A58").Value = 45424.58333333334
$ws.Range("B58").Value = 4360
$ws.Range("A59").Value = 45424.59375
$ws.Range("B59").Value = 4370
$ws.Range("A60").Value = 45424.60416666666
$ws.Range("B60").Value = 4390
$ws.Range("A61").Value = 45424.61458333334
$ws.Range("B61").Value = 4410
$ws.Range("A62").Value = 45424.625
$ws.Range("B62").Value = 4440
$ws.Range("A63").Value = 45424.63541666666
$ws.Range("B63").Value = 4490
$ws.Range("A64").Value = 45424.64583333334
$ws.Range("B64").Value = 4550
$ws.Range("A65").Value = 45424.65625
$ws.Range("B65").Value = 4620
$ws.Range("A66").Value = 45424.66666666666
$ws.Range("B66").Value = 4690
$ws.Range("A67").Value = 45424.67708333334
$ws.Range("B67").Value = 4770
$ws.Range("A68").Value = 45424.6875
$ws.Range("B68").Value = 4850
$ws.Range("A69").Value = 45424.69791666666
$ws.Range("B69").Value = 4930
$ws.Range("A70").Value = 45424.70833333334
$ws.Range("B70").Value = 5010
$ws.Range("A71").Value = 45424.71875
$ws.Range("B71").Value = 5090
$ws.Range("A72").Value = 45424.72916666666
$ws.Range("B72").Value = 5190
$ws.Range("A73").Value = 45424.73958333334
$ws.Range("B73").Value = 5290
$ws.Range("A74").Value = 45424.75
$ws.Range("B74").Value = 5400
$ws.Range("A75").Value = 45424.76041666666
$ws.Range("B75").Value = 5510
$ws.Range("A76").Value = 45424.77083333334
$ws.Range("B76").Value = 5630
$ws.Range("A77").Value = 45424.78125
$ws.Range("B77").Value = 5740
$ws.Range("A78").Value = 45424.79166666666
$ws.Range("B78").Value = 5860
$ws.Range("A79").Value = 45424.80208333334
$ws.Range("B79").Value = 5970
$ws.Range("A80").Value = 45424.8125
$ws.Range("B80").Value = 6040
$ws.Range("A81").Value = 45424.82291666666
$ws.Range("B81").Value = 6110
$ws.Range("A82").Value = 45424.83333333334
$ws.Range("B82").Value = 6120
$ws.Range("A83").Value = 45424.84375
$ws.Range("B83").Value = 6110
$ws.Range("A84").Value = 45424.85416666666
$ws.Range("B84").Value = 6080
$ws.Range("A85").Value = 45424.86458333334
$ws.Range("B85").Value = 6000
$ws.Range("A86").Value = 45424.875
$ws.Range("B86").Value = 5870
$ws.Range("A87").Value = 45424.88541666666
$ws.Range("B87").Value = 5750
$ws.Range("A88").Value = 45424.89583333334
$ws.Range("B88").Value = 5620
$ws.Range("A89").Value = 45424.90625
$ws.Range("B89").Value = 5490
$ws.Range("A90").Value = 45424.91666666666
$ws.Range("B90").Value = 5360
$ws.Range("A91").Value = 45424.92708333334
$ws.Range("B91").Value = 5240
$ws.Range("A92").Value = 45424.9375
$ws.Range("B92").Value = 5130
$ws.Range("A93").Value = 45424.94791666666
$ws.Range("B93").Value = 5020

# Add new rows 94-98, matching the date/style formatting of column A
$ws.Range("A94").Value = 45424.95833333334
$ws.Range("A94").NumberFormat = $ws.Range("A93").NumberFormat
$ws.Range("B94").Value = 5050
$ws.Range("A95").Value = 45424.96875
$ws.Range("A95").NumberFormat = $ws.Range("A93").NumberFormat
$ws.Range("B95").Value = 4990
$ws.Range("A96").Value = 45424.97916666666
$ws.Range("A96").NumberFormat = $ws.Range("A93").NumberFormat
$ws.Range("B96").Value = 4930
$ws.Range("A97").Value = 45424.98958333334
$ws.Range("A97").NumberFormat = $ws.Range("A93").NumberFormat
$ws.Range("B97").Value = 4880
$ws.Range("A98").Value = 45425
$ws.Range("A98").NumberFormat = $ws.Range("A93").NumberFormat
$ws.Range("B98").Value = 4840
